$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.Value = "'" + $val
    $rng.Style = $origStyle
}

Set-TextCell $ws 'D2' '266.21'
Set-TextCell $ws 'D3' '21.51'
Set-TextCell $ws 'D4' '6.116'
Set-TextCell $ws 'D5' '0.06106'
Set-TextCell $ws 'D6' '3.565'
Set-TextCell $ws 'D8' '1.354'
Set-TextCell $ws 'D9' '0.8196'
Set-TextCell $ws 'D10' '0.01342'
Set-TextCell $ws 'D11' '0.1580'
Set-TextCell $ws 'D12' '0.08082'
Set-TextCell $ws 'D13' '0.03443'
Set-TextCell $ws 'D14' '0.03206'
Set-TextCell $ws 'D15' '0.09227'
Set-TextCell $ws 'D16' '3.746'
Set-TextCell $ws 'D17' '0.001626'
Set-TextCell $ws 'D18' '0.04652'
Set-TextCell $ws 'D19' '0.006316'
Set-TextCell $ws 'D20' '0.006147'
Set-TextCell $ws 'D21' '0.001068'
Set-TextCell $ws 'D22' '0.0001500'
Set-TextCell $ws 'D24' '2.256'
Set-TextCell $ws 'D26' '0.1242'
Set-TextCell $ws 'D28' '0.0002713'
Set-TextCell $ws 'D40' '0.04591'
Set-TextCell $ws 'D41' '0.006992'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextCell $ws 'D42' '0.004000'
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextCell $ws 'D43' '0.1115'
$ws.Range('E43').Value = '42BKEXTokenBKK'
Set-TextCell $ws 'D45' '0.00005830'
Set-TextCell $ws 'D46' '0.0009900'
Set-TextCell $ws 'D47' '0.00000000750'
Set-TextCell $ws 'D48' '0.8025'
Set-TextCell $ws 'D50' '0.00001900'
Set-TextCell $ws 'D51' '0.01240'
